# Update crypto price/volume figures per the latest refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.698.05"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "1.961.04"
$ws.Range("E3").Value = "  +0.83%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.93"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.374"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("E10").Value = "  -3.93%  "

$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.08%  "

$ws.Range("D13").Value = "2.248.73"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.825"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").Value = "1.960.83"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").Value = "36.573.20"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").Value = "0.0₃0861"
$ws.Range("E20").Value = "  -1.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("E22").Value = "  +1.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.48%  "

$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.141"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("E30").Value = "  +0.88%  "

$ws.Range("E31").Value = "  -2.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0617"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.61%  "

$ws.Range("E34").Value = "  -0.58%  "

$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.23%  "

$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("E40").Value = "  +3.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.13%  "

$ws.Range("E42").Value = "  +1.03%  "

$ws.Range("E43").Value = "  -1.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").Value = "1.357.54"
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.24%  "

$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").Value = "2.139.19"
$ws.Range("E50").Value = "  +1.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.76%  "
